$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generische Lebensmittel")

# Decrement the ID values in column A for rows 2 through 238 (so they start at 0 instead of 1)
for ($r = 2; $r -le 238; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    $cell.Value = $current - 1
}

# Update the frozen-pane view to show the top of the sheet with A2 selected
$ws.Activate()
$ws.Range("A2").Select()
